# Prima prova di introduzione della funzione obiettivo in una ricerca locale
# (insert intramacchina - LS2)
#
# Adds a new column S ("funzione obiettivo" contribution per row) to the
# "Schedulazione" sheet, populating rows 2-27 with the values from the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(7, 7, 7, 7, 7, 7, 7, 7, 4, 2, 1, 1, 1, 1, 1, 1, 1, 1, 1, 2, 2, 2, 7, 4, 4, 2)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 19).Value = $values[$i]
}
